# Updates to Data and Affects
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix abbreviated "w" -> "with" in the description column (rows are Table1 data rows)
$ws.Range("B3").Value  = "Beaming Face with Smiling Eyes"
$ws.Range("B5").Value  = "Grinning Face with Big Eyes"
$ws.Range("B6").Value  = "Grinning Face with Smiling Eyes"
$ws.Range("B8").Value  = "Smiling Face with Halo"
$ws.Range("B13").Value = "Downcast Face with Sweat"
$ws.Range("B18").Value = "Face with Steam From Nose"
$ws.Range("B21").Value = "Anxious Face with Sweat"

# Update the view state: scroll the frozen pane so column B is the left-most
# visible column, and move the selection to B32 (below the table).
$win = $excel.ActiveWindow
$win.ScrollColumn = 2
$win.ScrollRow = 1
$ws.Range("B32").Select()
